# Update "想去人数" (interest count) figures to match latest scrape (commit 456a3b4)

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - rows 4-9, column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 70
$wsExpo.Range("F5").Value = 5149
$wsExpo.Range("F6").Value = 174
$wsExpo.Range("F7").Value = 51
$wsExpo.Range("F8").Value = 96
$wsExpo.Range("F9").Value = 333

# Sheet "全部类型" (all types) - rows 8-12 and 14, column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 70
$wsAll.Range("F9").Value = 5149
$wsAll.Range("F10").Value = 174
$wsAll.Range("F11").Value = 51
$wsAll.Range("F12").Value = 96
$wsAll.Range("F14").Value = 333
